$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.756.61'
$ws.Range('E2').Value = '  +2.35%  '

$ws.Range('D3').Value = '3.566.01'
$ws.Range('E3').Value = '  +1.33%  '

$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').Value = '584.13'
$ws.Range('E5').Value = '  +2.31%  '

$ws.Range('D6').Value = '185.56'
$ws.Range('E6').Value = '  +1.53%  '

$ws.Range('D7').Value = '''0.635'
$ws.Range('E7').Value = '  +3.34%  '

$ws.Range('D8').Value = '3.554.68'
$ws.Range('E8').Value = '  +1.18%  '

$ws.Range('E9').Value = '  +0.00%  '

$ws.Range('D10').Value = '''0.220'
$ws.Range('E10').Value = '  +18.93%  '

$ws.Range('E11').Value = '  +2.69%  '

$ws.Range('D12').Value = '''54.54'
$ws.Range('E12').Value = '  +1.33%  '

$ws.Range('D13').Value = '''0.0000318'

$ws.Range('E14').Value = '  +0.46%  '

$ws.Range('D15').Value = '4.133.74'
$ws.Range('E15').Value = '  +1.76%  '

$ws.Range('D16').Value = '70.839.89'
$ws.Range('E16').Value = '  +2.55%  '

$ws.Range('D17').Value = '''19.31'
$ws.Range('E17').Value = '  +0.35%  '

$ws.Range('D18').Value = '3.551.84'
$ws.Range('E18').Value = '  +1.16%  '

$ws.Range('D19').Value = '''576.03'
$ws.Range('E19').Value = '  +7.04%  '

$ws.Range('D20').Value = '12.44'
$ws.Range('E20').Value = '  -0.35%  '

$ws.Range('E21').Value = '  +0.79%  '

$ws.Range('D22').Value = '1.01'
$ws.Range('E22').Value = '  -2.12%  '

$ws.Range('D23').Value = '''17.72'
$ws.Range('E23').Value = '  -13.91%  '

$ws.Range('D24').Value = '5.04'
$ws.Range('E24').Value = '  +0.41%  '

$ws.Range('D25').Value = '''4.57'
$ws.Range('E25').Value = '  +4.24%  '

$ws.Range('D26').Value = '''95.62'
$ws.Range('E26').Value = '  +1.24%  '

$ws.Range('D27').Value = '11.26'
$ws.Range('E27').Value = '  +2.66%  '

$ws.Range('E28').Value = '  +1.41%  '

$ws.Range('D29').Value = '''9.14'
$ws.Range('E29').Value = '  +0.13%  '

$ws.Range('D30').Value = '''32.52'
$ws.Range('E30').Value = '  +3.26%  '

$ws.Range('D31').Value = '7.25'
$ws.Range('E31').Value = '  -0.28%  '

$ws.Range('D32').Value = '''12.29'
$ws.Range('E32').Value = '  -2.98%  '

$ws.Range('D33').Value = '''0.117'
$ws.Range('E33').Value = '  +2.61%  '

$ws.Range('D34').Value = '''63.49'
$ws.Range('E34').Value = '  -0.61%  '

$ws.Range('D35').Value = '''3.39'
$ws.Range('E35').Value = '  +10.25%  '

$ws.Range('D36').Value = '549.25'
$ws.Range('E36').Value = '  -4.45%  '

$ws.Range('D37').Value = '''0.415'
$ws.Range('E37').Value = '  +4.22%  '

$ws.Range('B38').Value = 'InjectiveProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D38').Value = '''37.90'
$ws.Range('E38').Value = '  -0.90%  '

$ws.Range('B39').Value = 'Dai'
$ws.Range('C39').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D39').Value = '''1.00'
$ws.Range('E39').Value = '  +0.10%  '

$ws.Range('B40').Value = 'dogwifhat'
$ws.Range('C40').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D40').Value = '''3.36'
$ws.Range('E40').Value = '  +6.51%  '

$ws.Range('B41').Value = 'PEPE'
$ws.Range('C41').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D41').Value = '0.0₃0794'
$ws.Range('E41').Value = '  +4.28%  '

$ws.Range('D42').Value = '3.566.38'
$ws.Range('E42').Value = '  +11.92%  '

$ws.Range('E43').Value = '  +2.45%  '

$ws.Range('E44').Value = '  +2.55%  '

$ws.Range('D45').Value = '''3.53'
$ws.Range('E45').Value = '  +0.20%  '

$ws.Range('D46').Value = '''0.0446'
$ws.Range('E46').Value = '  +1.12%  '

$ws.Range('E47').Value = '  -1.19%  '

$ws.Range('D48').Value = '''9.33'
$ws.Range('E48').Value = '  +1.90%  '

$ws.Range('D49').Value = '''0.139'
$ws.Range('E49').Value = '  +3.57%  '

$ws.Range('B50').Value = 'FLOKI'
$ws.Range('C50').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D50').Value = '''0.000265'
$ws.Range('E50').Value = '  +17.65%  '

$ws.Range('E51').Value = '  +0.02%  '
